# Add 2022-Q4 data:
#  1. Insert a new "2022-Q4" worksheet (cloned from the existing "2022-Q3"
#     sheet so it inherits identical structure/styles), positioned right
#     after "总计" and before "2022-Q3".
#  2. Trim it down to the 12 data rows (+ header) that the new quarter needs
#     and overwrite its contents with the 2022-Q4 fund holdings.
#  3. Insert a new row 2 at the top of the "总计" summary sheet for the new
#     quarter and renumber the running index column.

$wb = $excel.ActiveWorkbook

# ---- 1. Clone "2022-Q3" (sheet index 2) to create "2022-Q4" ----------------
$template = $wb.Worksheets.Item(2)
$template.Copy($template, $null)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# The template has 21 data rows (rows 2-22); the new quarter only has 12
# data rows (rows 2-13), so drop the extra rows.
$newSheet.Rows("14:22").Delete()

# ---- 2. Fill in the 2022-Q4 fund holdings ----------------------------------
$rows = @(
    @("001735", "广发百发大数据策略成长灵活配置混合E", "8.88", "90.16", "1.44", "0.1279", 4),
    @("001734", "广发百发大数据策略成长灵活配置混合A", "4.29", "90.16", "1.44", "0.0618", 4),
    @("007832", "博道伍佰智航股票C", "6.00", "93.04", "0.88", "0.0528", 9),
    @("014135", "中欧金安量化混合A", "7.09", "90.01", "0.63", "0.0447", 8),
    @("001990", "中欧数据挖掘多因子灵活配置混合A", "3.26", "90.14", "0.75", "0.0244", 3),
    @("007831", "博道伍佰智航股票A", "2.75", "93.04", "0.88", "0.0242", 9),
    @("011410", "中信建投量化进取6个月持有期混合A", "5.05", "69.90", "0.46", "0.0232", 5),
    @("002137", "诺安利鑫灵活配置混合A", "0.44", "89.87", "4.96", "0.0218", 1),
    @("004234", "中欧数据挖掘多因子灵活配置混合C", "1.93", "90.14", "0.75", "0.0145", 3),
    @("011411", "中信建投量化进取6个月持有期混合C", "1.56", "69.90", "0.46", "0.0072", 5),
    @("014136", "中欧金安量化混合C", "1.07", "90.01", "0.63", "0.0067", 8),
    @("014521", "诺安利鑫灵活配置混合C", "0.01", "89.87", "4.96", "0.0005", 1)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $item = $rows[$i]

    $newSheet.Cells.Item($r, 1).Value = $i

    # Fund code looks numeric ("001735") - force text so the leading zero
    # survives, then drop the quote-prefix styling it picks up.
    $codeCell = $newSheet.Cells.Item($r, 2)
    $codeCell.Value = "'" + $item[0]
    $codeCell.Style = "Normal"

    $newSheet.Cells.Item($r, 3).Value = $item[1]

    # Columns D-G hold numeric-looking figures that are stored as text in
    # the source data; force text the same way as the fund code.
    for ($c = 4; $c -le 7; $c++) {
        $cell = $newSheet.Cells.Item($r, $c)
        $cell.Value = "'" + $item[$c - 2]
        $cell.Style = "Normal"
    }

    $newSheet.Cells.Item($r, 8).Value = $item[6]
}

# ---- 3. Add the new quarter to the "总计" summary sheet --------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows("2:2").Insert()

# Restore column-A styling (the Insert() doesn't copy it across) by cloning
# the format from the row below, then reset the rest of the new row to the
# default (unstyled) look used by every other data row.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").Style = "Normal"

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 12
$summary.Cells.Item(2, 4).Value = 0.41

# Renumber the running index in column A for the rows that got pushed down.
for ($r = 3; $r -le 7; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# Leave the original first sheet active, like before the edit.
$summary.Activate()
